# Lab Exam 03 grading workbook - fill in "Total Points" (column E) scores
# to match the "Points for grading" (column D) for the first two rubric
# sections (Constructor/Getter/toString/Compilation questions and the
# Customer Class questions), matching the grader's marks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section 1 (Generic Class questions, rows 3-6): award full marks in column E
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Section 2 (Customer Class questions, rows 10-14): award full marks in column E
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Recalculate so the section subtotals (E7, E15) and grand total (E38) update
$excel.CalculateFull()

# Move selection/viewport back to the top and select E15 (the subtotal the
# grader just finished computing), clearing the previous scroll position.
$ws.Range("A1").Select() | Out-Null
$ws.Range("E15").Select() | Out-Null
